$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at the bottom of the table, copying format from the row above
# (xlShiftDown, xlFormatFromLeftOrAbove) so the new rows pick up the same per-column
# styles (e.g. the email column's style and the is_active boolean column's style)
# already used throughout the table.
$ws.Rows.Item(31).Insert(-4121, 0)
$ws.Rows.Item(32).Insert(-4121, 0)

# Populate row 32 (John Doe) first so "John Doe" / "john.doe@xyz.com" land earlier
# in the shared-strings table than "Jane Smith" / "jane.smith@xyz.com".
$ws.Cells.Item(32,1).Value = 110031
$ws.Cells.Item(32,2).Value = 9317596767
$ws.Cells.Item(32,3).Value = "John Doe"
$ws.Cells.Item(32,4).Value = "john.doe@xyz.com"
$ws.Cells.Item(32,5).Value = 818876431
$ws.Cells.Item(32,6).Value = "ACT"
$ws.Cells.Item(32,7).Value = "eng"
$ws.Cells.Item(32,8).Value = "PWD"
$ws.Cells.Item(32,9).Value = $true
$ws.Cells.Item(32,10).Value = "superadmin"
$ws.Cells.Item(32,11).Value = "now()"
$ws.Cells.Item(32,12).Value = "now()"

# Row 31 (Jane Smith)
$ws.Cells.Item(31,1).Value = 110030
$ws.Cells.Item(31,2).Value = 9317596768
$ws.Cells.Item(31,3).Value = "Jane Smith"
$ws.Cells.Item(31,4).Value = "jane.smith@xyz.com"
$ws.Cells.Item(31,5).Value = 818876432
$ws.Cells.Item(31,6).Value = "ACT"
$ws.Cells.Item(31,7).Value = "eng"
$ws.Cells.Item(31,8).Value = "PWD"
$ws.Cells.Item(31,9).Value = $true
$ws.Cells.Item(31,10).Value = "superadmin"
$ws.Cells.Item(31,11).Value = "now()"
$ws.Cells.Item(31,12).Value = "now()"

$ws.Range("F30").Select()
